# Auto-generated edit script applying the diff's cell-value updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 569.81396
$ws.Range("J33").Value = 207
$ws.Range("L33").Value = 207
$ws.Range("N33").Value = -665
$ws.Range("H112").Value = 812614.4
$ws.Range("I112").Value = 654
$ws.Range("J112").Value = 989127.5
$ws.Range("K112").Value = 1962
$ws.Range("L112").Value = 2967382.5
$ws.Range("M112").Value = -854
$ws.Range("N112").Value = -2969598.5
$ws.Range("H127").Value = 890.6799999999999
$ws.Range("I127").Value = 542.4
$ws.Range("J127").Value = 1413.1
$ws.Range("K127").Value = 1627.2
$ws.Range("L127").Value = 4239.299999999999
$ws.Range("M127").Value = 3332.8
$ws.Range("N127").Value = -14159.3
$ws.Range("H129").Value = 7006.8335
$ws.Range("I129").Value = 372.125
$ws.Range("J129").Value = 8902.464
$ws.Range("K129").Value = 1116.375
$ws.Range("L129").Value = 26707.392
$ws.Range("M129").Value = 3883.625
$ws.Range("N129").Value = -36707.392
$ws.Range("H132").Value = 2903.5513
$ws.Range("I132").Value = 2626.0442
$ws.Range("J132").Value = 4790.6
$ws.Range("K132").Value = 7878.132599999999
$ws.Range("L132").Value = 14371.8
$ws.Range("M132").Value = -5348.132599999999
$ws.Range("N132").Value = -19431.8
$ws.Range("H138").Value = 3436.8225
$ws.Range("I138").Value = 1943.5883
$ws.Range("J138").Value = 4000.9333
$ws.Range("K138").Value = 5830.7649
$ws.Range("L138").Value = 12002.7999
$ws.Range("M138").Value = -690.7649000000001
$ws.Range("N138").Value = -22282.7999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 142973150
$ws.Range("J88").Value = 166801520
$ws.Range("L88").Value = 166801520
$ws.Range("N88").Value = -166802332
$ws.Range("H91").Value = 142973150
$ws.Range("J91").Value = 166801520
$ws.Range("L91").Value = 166801520
$ws.Range("N91").Value = -166804328
$ws.Range("H97").Value = 3730.0833
$ws.Range("I97").Value = 4638.8887
$ws.Range("J97").Value = 1003.6667
$ws.Range("K97").Value = 4638.8887
$ws.Range("L97").Value = 1003.6667
$ws.Range("M97").Value = -4142.8887
$ws.Range("N97").Value = -1995.6667
$ws.Range("H105").Value = 35250
$ws.Range("J105").Value = 35250
$ws.Range("L105").Value = 35250
$ws.Range("N105").Value = -42238
$ws.Range("H132").Value = 18761.322
$ws.Range("I132").Value = 28523.19
$ws.Range("J132").Value = 2343.6365
$ws.Range("K132").Value = 85569.56999999999
$ws.Range("L132").Value = 7030.9095
$ws.Range("M132").Value = -83039.56999999999
$ws.Range("N132").Value = -12090.9095

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 28747.273
$ws.Range("I94").Value = 26097
$ws.Range("J94").Value = 30261.715
$ws.Range("K94").Value = 26097
$ws.Range("L94").Value = 30261.715
$ws.Range("M94").Value = -25646
$ws.Range("N94").Value = -31163.715
$ws.Range("H99").Value = 1743.5
$ws.Range("I99").Value = 1485.7142
$ws.Range("K99").Value = 1485.7142
$ws.Range("M99").Value = 12.28580000000011
$ws.Range("H107").Value = 7727.625
$ws.Range("I107").Value = 11466.2
$ws.Range("J107").Value = 1496.6666
$ws.Range("K107").Value = 11466.2
$ws.Range("L107").Value = 1496.6666
$ws.Range("M107").Value = -9546.200000000001
$ws.Range("N107").Value = -5336.6666
$ws.Range("H110").Value = 35000
$ws.Range("J110").Value = 35000
$ws.Range("L110").Value = 35000
$ws.Range("N110").Value = -43180

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()
$ws.Range("H58").Value = 2191.9524
$ws.Range("I58").Value = 1544.7273
$ws.Range("J58").Value = 2903.9
$ws.Range("K58").Value = 1544.7273
$ws.Range("L58").Value = 2903.9
$ws.Range("M58").Value = -1341.7273
$ws.Range("N58").Value = -3309.9
$ws.Range("H99").Value = 4190.933
$ws.Range("I99").Value = 4707.1113
$ws.Range("J99").Value = 3416.6667
$ws.Range("K99").Value = 4707.1113
$ws.Range("L99").Value = 3416.6667
$ws.Range("M99").Value = -3209.1113
$ws.Range("N99").Value = -6412.6667
$ws.Range("H126").Value = 4190.933
$ws.Range("I126").Value = 4707.1113
$ws.Range("J126").Value = 3416.6667
$ws.Range("K126").Value = 14121.3339
$ws.Range("L126").Value = 10250.0001
$ws.Range("M126").Value = -11651.3339
$ws.Range("N126").Value = -15190.0001
$ws.Range("H134").Value = 3372.4524
$ws.Range("I134").Value = 3616.9678
$ws.Range("J134").Value = 2683.3635
$ws.Range("K134").Value = 10850.9034
$ws.Range("L134").Value = 8050.0905
$ws.Range("M134").Value = -8315.903399999999
$ws.Range("N134").Value = -13120.0905
$ws.Range("H136").Value = 2191.9524
$ws.Range("I136").Value = 1544.7273
$ws.Range("J136").Value = 2903.9
$ws.Range("K136").Value = 4634.1819
$ws.Range("L136").Value = 8711.700000000001
$ws.Range("M136").Value = -2084.1819
$ws.Range("N136").Value = -13811.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 772.4792
$ws.Range("I5").Value = 451.1154
$ws.Range("J5").Value = 1152.2727
$ws.Range("K5").Value = 1353.3462
$ws.Range("L5").Value = 3456.8181
$ws.Range("M5").Value = -1241.3462
$ws.Range("N5").Value = -3680.8181
$ws.Range("H131").Value = 1588044.6
$ws.Range("I131").Value = 4167030.5
$ws.Range("J131").Value = 976.5
$ws.Range("K131").Value = 12501091.5
$ws.Range("L131").Value = 2929.5
$ws.Range("M131").Value = -12496051.5
$ws.Range("N131").Value = -13009.5
$ws.Range("H135").Value = 772.4792
$ws.Range("I135").Value = 451.1154
$ws.Range("J135").Value = 1152.2727
$ws.Range("K135").Value = 4060.0386
$ws.Range("L135").Value = 10370.4543
$ws.Range("M135").Value = -1525.0386
$ws.Range("N135").Value = -15440.4543
$ws.Range("H136").Value = 1203.4706
$ws.Range("I136").Value = 1208.4286
$ws.Range("K136").Value = 3625.2858
$ws.Range("M136").Value = 1474.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3919.0625
$ws.Range("I80").Value = 5584.1665
$ws.Range("J80").Value = 2920
$ws.Range("K80").Value = 5584.1665
$ws.Range("L80").Value = 2920
$ws.Range("M80").Value = -4586.1665
$ws.Range("N80").Value = -4916
$ws.Range("H83").Value = 3919.0625
$ws.Range("I83").Value = 5584.1665
$ws.Range("J83").Value = 2920
$ws.Range("K83").Value = 27920.8325
$ws.Range("L83").Value = 14600
$ws.Range("M83").Value = -22928.8325
$ws.Range("N83").Value = -24584
$ws.Range("H99").Value = 4940.625
$ws.Range("I99").Value = 2789.2856
$ws.Range("J99").Value = 20000
$ws.Range("K99").Value = 2789.2856
$ws.Range("L99").Value = 20000
$ws.Range("M99").Value = -543.2856000000002
$ws.Range("N99").Value = -24492
$ws.Range("H126").Value = 2637.1333
$ws.Range("I126").Value = 1820.6666
$ws.Range("J126").Value = 3453.6
$ws.Range("K126").Value = 5461.9998
$ws.Range("L126").Value = 10360.8
$ws.Range("M126").Value = -2991.9998
$ws.Range("N126").Value = -15300.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1416.52
$ws.Range("I40").Value = 1240.65
$ws.Range("J40").Value = 2120
$ws.Range("K40").Value = 1240.65
$ws.Range("L40").Value = 2120
$ws.Range("M40").Value = -1104.65
$ws.Range("N40").Value = -2392
$ws.Range("H93").Value = 3275
$ws.Range("I93").Value = 3275
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 3275
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -2027
$ws.Range("N93").ClearContents()
$ws.Range("H122").Value = 4265.343
$ws.Range("I122").Value = 4190.2188
$ws.Range("J122").Value = 5066.6665
$ws.Range("K122").Value = 12570.6564
$ws.Range("L122").Value = 15199.9995
$ws.Range("M122").Value = -10120.6564
$ws.Range("N122").Value = -20099.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 1585.75
$ws.Range("I8").Value = 103
$ws.Range("J8").Value = 2080
$ws.Range("K8").Value = 103
$ws.Range("L8").Value = 2080
$ws.Range("M8").Value = 37
$ws.Range("N8").Value = -2360
$ws.Range("H126").Value = 3288.25
$ws.Range("I126").Value = 5897
$ws.Range("J126").Value = 679.5
$ws.Range("K126").Value = 17691
$ws.Range("L126").Value = 2038.5
$ws.Range("M126").Value = -15221
$ws.Range("N126").Value = -6978.5
